# Daily attendance processing - reverse order of comma-separated
# "Recorded By" entries in column G of the session analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $reversed = @()
        for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
            $reversed += $trimmed[$i]
        }

        $newVal = [string]::Join(", ", $reversed)
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
